$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("E2").Value = "iaest-measure:provincia"
$ws.Range("E3").Value = "medida"
$ws.Range("D4").Value = "URI-Comunidad"
$ws.Range("E4").Value = "xsd:int"

$ws.Rows.Item(5).Delete()
